$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was collected (2022-12-02). Insert two new
# data rows right below the header block of existing rows for this
# product (rows 26:27), pushing the rest of the table down by two rows.
$ws.Rows("26:27").Insert()

$fecha = Get-Date -Year 2022 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0

# Row 26: Cebollín, Primera
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C26").Value = 'Ñuble'
$ws.Range("D26").Value = $fecha
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112037
$ws.Range("G26").Value = 'Cebollín'
$ws.Range("H26").Value = 'Sin especificar'
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 600
$ws.Range("L26").Value = 700
$ws.Range("M26").Value = 650
$ws.Range("N26").Value = '$/paquete 6 unidades'
$ws.Range("O26").Value = 'Provincia de Diguillín'
$ws.Range("P26").Value = 108
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = 'Hortaliza'

# Row 27: Cebollín, Segunda
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C27").Value = 'Ñuble'
$ws.Range("D27").Value = $fecha
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112037
$ws.Range("G27").Value = 'Cebollín'
$ws.Range("H27").Value = 'Sin especificar'
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = 500
$ws.Range("N27").Value = '$/paquete 6 unidades'
$ws.Range("O27").Value = 'Provincia de Diguillín'
$ws.Range("P27").Value = 83
$ws.Range("Q27").Value = 6
$ws.Range("R27").Value = 'Hortaliza'
